# Auto-generated Excel COM-interop script to apply market-price/profit data refresh
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6997.909
$ws.Range("I40").Value = 5999.4
$ws.Range("J40").Value = 7830
$ws.Range("K40").Value = 5999.4
$ws.Range("L40").Value = 7830
$ws.Range("M40").Value = -5824.4
$ws.Range("N40").Value = -8180
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H100").Value = 7999.25
$ws.Range("I100").Value = 7997
$ws.Range("K100").Value = 7997
$ws.Range("M100").Value = -7456
$ws.Range("H113").Value = 1540
$ws.Range("I113").Value = 1540
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1540
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1714
$ws.Range("N113").ClearContents()
$ws.Range("H138").Value = 1821.2727
$ws.Range("I138").Value = 1670.4445
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 5011.333500000001
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = 128.6664999999994
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5339.839
$ws.Range("I32").Value = 5339.839
$ws.Range("K32").Value = 5339.839
$ws.Range("M32").Value = -5052.839
$ws.Range("H63").Value = 1830.8
$ws.Range("I63").Value = 1788.5
$ws.Range("K63").Value = 1788.5
$ws.Range("M63").Value = -1102.5
$ws.Range("H66").Value = 1830.8
$ws.Range("I66").Value = 1788.5
$ws.Range("K66").Value = 8942.5
$ws.Range("M66").Value = -5510.5
$ws.Range("H74").Value = 838.5625
$ws.Range("I74").Value = 854.4666999999999
$ws.Range("K74").Value = 854.4666999999999
$ws.Range("M74").Value = 19.53330000000005
$ws.Range("H77").Value = 838.5625
$ws.Range("I77").Value = 854.4666999999999
$ws.Range("K77").Value = 4272.3335
$ws.Range("M77").Value = 95.66650000000027
$ws.Range("H88").Value = 2527.2222
$ws.Range("I88").Value = 4800
$ws.Range("J88").Value = 1877.8572
$ws.Range("K88").Value = 4800
$ws.Range("L88").Value = 1877.8572
$ws.Range("M88").Value = -4394
$ws.Range("N88").Value = -2689.8572
$ws.Range("H91").Value = 2527.2222
$ws.Range("I91").Value = 4800
$ws.Range("J91").Value = 1877.8572
$ws.Range("K91").Value = 4800
$ws.Range("L91").Value = 1877.8572
$ws.Range("M91").Value = -3396
$ws.Range("N91").Value = -4685.8572
$ws.Range("H97").Value = 630
$ws.Range("I97").Value = 604.63635
$ws.Range("K97").Value = 604.63635
$ws.Range("M97").Value = -108.63635
$ws.Range("H102").Value = 6765.875
$ws.Range("I102").Value = 3463.1428
$ws.Range("K102").Value = 3463.1428
$ws.Range("M102").Value = -1841.1428
$ws.Range("H113").Value = 19000
$ws.Range("J113").Value = 19000
$ws.Range("L113").Value = 19000
$ws.Range("N113").Value = -27678
$ws.Range("H122").Value = 1833.2858
$ws.Range("I122").Value = 1833.2858
$ws.Range("K122").Value = 5499.857400000001
$ws.Range("M122").Value = -3049.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7419.9
$ws.Range("J86").Value = 7524.875
$ws.Range("L86").Value = 7524.875
$ws.Range("N86").Value = -9770.875
$ws.Range("H89").Value = 7419.9
$ws.Range("J89").Value = 7524.875
$ws.Range("L89").Value = 37624.375
$ws.Range("N89").Value = -48856.375
$ws.Range("H105").Value = 1714.5
$ws.Range("I105").Value = 2030
$ws.Range("J105").Value = 1399
$ws.Range("K105").Value = 2030
$ws.Range("L105").Value = 1399
$ws.Range("M105").Value = -283
$ws.Range("N105").Value = -4893
$ws.Range("H140").Value = 143556
$ws.Range("J140").Value = 139926.67
$ws.Range("L140").Value = 139926.67
$ws.Range("N140").Value = -150286.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5674.282
$ws.Range("I31").Value = 2101.15
$ws.Range("J31").Value = 9435.474
$ws.Range("K31").Value = 2101.15
$ws.Range("L31").Value = 9435.474
$ws.Range("M31").Value = -1806.15
$ws.Range("N31").Value = -10025.474
$ws.Range("H34").Value = 5674.282
$ws.Range("I34").Value = 2101.15
$ws.Range("J34").Value = 9435.474
$ws.Range("K34").Value = 2101.15
$ws.Range("L34").Value = 9435.474
$ws.Range("M34").Value = -1899.15
$ws.Range("N34").Value = -9839.474
$ws.Range("H82").Value = 50181
$ws.Range("J82").Value = 50181
$ws.Range("L82").Value = 50181
$ws.Range("N82").Value = -50903
$ws.Range("H85").Value = 50181
$ws.Range("J85").Value = 50181
$ws.Range("L85").Value = 50181
$ws.Range("N85").Value = -52677
$ws.Range("H105").Value = 1576.4286
$ws.Range("I105").Value = 1263.6666
$ws.Range("K105").Value = 1263.6666
$ws.Range("M105").Value = 483.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 388.83334
$ws.Range("I44").Value = 266.6
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 799.8000000000001
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = -401.8000000000001
$ws.Range("N44").Value = -3796
$ws.Range("H54").Value = 3582.25
$ws.Range("I54").Value = 3499.8
$ws.Range("J54").Value = 3994.5
$ws.Range("K54").Value = 10499.4
$ws.Range("L54").Value = 11983.5
$ws.Range("M54").Value = -9940.400000000001
$ws.Range("N54").Value = -13101.5
$ws.Range("H60").Value = 678.5217
$ws.Range("I60").Value = 300
$ws.Range("J60").Value = 1091.4546
$ws.Range("K60").Value = 900
$ws.Range("L60").Value = 3274.3638
$ws.Range("M60").Value = -649
$ws.Range("N60").Value = -3776.3638
$ws.Range("H86").Value = 382.2
$ws.Range("I86").Value = 260.25
$ws.Range("J86").Value = 463.5
$ws.Range("K86").Value = 780.75
$ws.Range("L86").Value = 1390.5
$ws.Range("M86").Value = 405.25
$ws.Range("N86").Value = -3762.5
$ws.Range("H89").Value = 382.2
$ws.Range("I89").Value = 260.25
$ws.Range("J89").Value = 463.5
$ws.Range("K89").Value = 2342.25
$ws.Range("L89").Value = 4171.5
$ws.Range("M89").Value = 3585.75
$ws.Range("N89").Value = -16027.5
$ws.Range("H98").Value = 417.25
$ws.Range("I98").Value = 146.5
$ws.Range("K98").Value = 439.5
$ws.Range("M98").Value = 1058.5
$ws.Range("H128").Value = 1900000
$ws.Range("I128").Value = 1900000
$ws.Range("K128").Value = 5700000
$ws.Range("M128").Value = -5695020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5332.3335
$ws.Range("I70").Value = 5332.3335
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5332.3335
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5062.3335
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 5332.3335
$ws.Range("I73").Value = 5332.3335
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5332.3335
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4396.3335
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 4283.3335
$ws.Range("J80").Value = 6000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7996
$ws.Range("H83").Value = 4283.3335
$ws.Range("J83").Value = 6000
$ws.Range("L83").Value = 30000
$ws.Range("N83").Value = -39984
$ws.Range("H97").Value = 962.46155
$ws.Range("I97").Value = 910.63635
$ws.Range("K97").Value = 910.63635
$ws.Range("M97").Value = -414.63635
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 162215.14
$ws.Range("I132").Value = 274505.75
$ws.Range("K132").Value = 823517.25
$ws.Range("M132").Value = -820987.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1164.3
$ws.Range("I22").Value = 738.1111
$ws.Range("K22").Value = 738.1111
$ws.Range("M22").Value = -443.1111
$ws.Range("H27").Value = 1164.3
$ws.Range("I27").Value = 738.1111
$ws.Range("K27").Value = 738.1111
$ws.Range("M27").Value = -631.1111
$ws.Range("H38").Value = 30343.666
$ws.Range("J38").Value = 30343.666
$ws.Range("L38").Value = 30343.666
$ws.Range("N38").Value = -31163.666
$ws.Range("H82").Value = 7987.5
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 7987.5
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H132").Value = 7340.857
$ws.Range("I132").Value = 3796.3333
$ws.Range("K132").Value = 11388.9999
$ws.Range("M132").Value = -8858.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3082.8333
$ws.Range("I2").Value = 3082.8333
$ws.Range("K2").Value = 3082.8333
$ws.Range("M2").Value = -2970.8333
$ws.Range("H15").Value = 57007
$ws.Range("J15").Value = 57007
$ws.Range("L15").Value = 57007
$ws.Range("N15").Value = -57583
$ws.Range("H62").Value = 8499.5
$ws.Range("J62").Value = 8726.727999999999
$ws.Range("L62").Value = 8726.727999999999
$ws.Range("N62").Value = -9974.727999999999
$ws.Range("H65").Value = 8499.5
$ws.Range("J65").Value = 8726.727999999999
$ws.Range("L65").Value = 43633.64
$ws.Range("N65").Value = -49873.64
